$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RG value in C2
$ws.Range("C2").Value = 333333333

# Update CELULAR column (K2:K6) with the new phone number text
$ws.Range("K2").Value = "14981225509"
$ws.Range("K3").Value = "14981225509"
$ws.Range("K4").Value = "14981225509"
$ws.Range("K5").Value = "14981225509"
$ws.Range("K6").Value = "14981225509"

# Move the active selection to K10, matching the final cursor position
$ws.Range("K10").Select()
